$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BD_Times")

$ws1.Cells.Item(440,1).Value = "Tombense"
$ws1.Cells.Item(440,2).Value = 1
$ws1.Cells.Item(440,3).Value = 0
$ws1.Cells.Item(440,4).Value = 0
$ws1.Cells.Item(440,5).Value = 0
$ws1.Cells.Item(440,6).Value = 0
$ws1.Cells.Item(440,7).Value = 0
$ws1.Cells.Item(440,8).Value = 2
$ws1.Cells.Item(440,9).Value = 5

$ws1.Cells.Item(441,1).Value = "Sport"
$ws1.Cells.Item(441,2).Value = 0
$ws1.Cells.Item(441,3).Value = 0
$ws1.Cells.Item(441,4).Value = 0
$ws1.Cells.Item(441,5).Value = 0
$ws1.Cells.Item(441,6).Value = 0
$ws1.Cells.Item(441,7).Value = 0
$ws1.Cells.Item(441,8).Value = 5
$ws1.Cells.Item(441,9).Value = 2

$ws1.Cells.Item(442,1).Value = "Criciuma"
$ws1.Cells.Item(442,2).Value = 1
$ws1.Cells.Item(442,3).Value = 1
$ws1.Cells.Item(442,4).Value = 0
$ws1.Cells.Item(442,5).Value = 0
$ws1.Cells.Item(442,6).Value = 2
$ws1.Cells.Item(442,7).Value = 0
$ws1.Cells.Item(442,8).Value = 7
$ws1.Cells.Item(442,9).Value = 3

$ws1.Cells.Item(443,1).Value = "Londrina"
$ws1.Cells.Item(443,2).Value = 0
$ws1.Cells.Item(443,3).Value = 0
$ws1.Cells.Item(443,4).Value = 1
$ws1.Cells.Item(443,5).Value = 0
$ws1.Cells.Item(443,6).Value = 0
$ws1.Cells.Item(443,7).Value = 2
$ws1.Cells.Item(443,8).Value = 3
$ws1.Cells.Item(443,9).Value = 7

$ws1.Cells.Item(444,1).Value = "Botafogo"
$ws1.Cells.Item(444,2).Value = 1
$ws1.Cells.Item(444,3).Value = 0
$ws1.Cells.Item(444,4).Value = 0
$ws1.Cells.Item(444,5).Value = 0
$ws1.Cells.Item(444,6).Value = 0
$ws1.Cells.Item(444,7).Value = 0
$ws1.Cells.Item(444,8).Value = 5
$ws1.Cells.Item(444,9).Value = 1

$ws1.Cells.Item(445,1).Value = "Ponte Preta"
$ws1.Cells.Item(445,2).Value = 0
$ws1.Cells.Item(445,3).Value = 0
$ws1.Cells.Item(445,4).Value = 0
$ws1.Cells.Item(445,5).Value = 0
$ws1.Cells.Item(445,6).Value = 0
$ws1.Cells.Item(445,7).Value = 0
$ws1.Cells.Item(445,8).Value = 1
$ws1.Cells.Item(445,9).Value = 5

$ws1.Cells.Item(446,1).Value = "Sampaio Correia"
$ws1.Cells.Item(446,2).Value = 1
$ws1.Cells.Item(446,3).Value = 1
$ws1.Cells.Item(446,4).Value = 1
$ws1.Cells.Item(446,5).Value = 1
$ws1.Cells.Item(446,6).Value = 1
$ws1.Cells.Item(446,7).Value = 1
$ws1.Cells.Item(446,8).Value = 7
$ws1.Cells.Item(446,9).Value = 7

$ws1.Cells.Item(447,1).Value = "CRB"
$ws1.Cells.Item(447,2).Value = 0
$ws1.Cells.Item(447,3).Value = 1
$ws1.Cells.Item(447,4).Value = 1
$ws1.Cells.Item(447,5).Value = 1
$ws1.Cells.Item(447,6).Value = 1
$ws1.Cells.Item(447,7).Value = 1
$ws1.Cells.Item(447,8).Value = 7
$ws1.Cells.Item(447,9).Value = 7

$ws1.Cells.Item(448,1).Value = "Vila Nova"
$ws1.Cells.Item(448,2).Value = 1
$ws1.Cells.Item(448,3).Value = 1
$ws1.Cells.Item(448,4).Value = 1
$ws1.Cells.Item(448,5).Value = 1
$ws1.Cells.Item(448,6).Value = 1
$ws1.Cells.Item(448,7).Value = 1
$ws1.Cells.Item(448,8).Value = 6
$ws1.Cells.Item(448,9).Value = 2

$ws1.Cells.Item(449,1).Value = "Avai"
$ws1.Cells.Item(449,2).Value = 0
$ws1.Cells.Item(449,3).Value = 1
$ws1.Cells.Item(449,4).Value = 1
$ws1.Cells.Item(449,5).Value = 1
$ws1.Cells.Item(449,6).Value = 1
$ws1.Cells.Item(449,7).Value = 1
$ws1.Cells.Item(449,8).Value = 2
$ws1.Cells.Item(449,9).Value = 6

$ws1.Cells.Item(450,1).Value = "Chapecoense"
$ws1.Cells.Item(450,2).Value = 1
$ws1.Cells.Item(450,3).Value = 0
$ws1.Cells.Item(450,4).Value = 1
$ws1.Cells.Item(450,5).Value = 0
$ws1.Cells.Item(450,6).Value = 0
$ws1.Cells.Item(450,7).Value = 1
$ws1.Cells.Item(450,8).Value = 4
$ws1.Cells.Item(450,9).Value = 4

$ws1.Cells.Item(451,1).Value = "Atletico-GO"
$ws1.Cells.Item(451,2).Value = 0
$ws1.Cells.Item(451,3).Value = 1
$ws1.Cells.Item(451,4).Value = 0
$ws1.Cells.Item(451,5).Value = 0
$ws1.Cells.Item(451,6).Value = 1
$ws1.Cells.Item(451,7).Value = 0
$ws1.Cells.Item(451,8).Value = 4
$ws1.Cells.Item(451,9).Value = 4

$ws1.Cells.Item(452,1).Value = "Guarani"
$ws1.Cells.Item(452,2).Value = 1
$ws1.Cells.Item(452,3).Value = 0
$ws1.Cells.Item(452,4).Value = 1
$ws1.Cells.Item(452,5).Value = 0
$ws1.Cells.Item(452,6).Value = 0
$ws1.Cells.Item(452,7).Value = 1
$ws1.Cells.Item(452,8).Value = 2
$ws1.Cells.Item(452,9).Value = 0

$ws1.Cells.Item(453,1).Value = "Juventude"
$ws1.Cells.Item(453,2).Value = 0
$ws1.Cells.Item(453,3).Value = 1
$ws1.Cells.Item(453,4).Value = 0
$ws1.Cells.Item(453,5).Value = 0
$ws1.Cells.Item(453,6).Value = 1
$ws1.Cells.Item(453,7).Value = 0
$ws1.Cells.Item(453,8).Value = 0
$ws1.Cells.Item(453,9).Value = 2

$ws1.Cells.Item(454,1).Value = "Vitoria"
$ws1.Cells.Item(454,2).Value = 1
$ws1.Cells.Item(454,3).Value = 1
$ws1.Cells.Item(454,4).Value = 0
$ws1.Cells.Item(454,5).Value = 0
$ws1.Cells.Item(454,6).Value = 1
$ws1.Cells.Item(454,7).Value = 0
$ws1.Cells.Item(454,8).Value = 4
$ws1.Cells.Item(454,9).Value = 9

$ws1.Cells.Item(455,1).Value = "Ceara"
$ws1.Cells.Item(455,2).Value = 0
$ws1.Cells.Item(455,3).Value = 0
$ws1.Cells.Item(455,4).Value = 1
$ws1.Cells.Item(455,5).Value = 0
$ws1.Cells.Item(455,6).Value = 0
$ws1.Cells.Item(455,7).Value = 1
$ws1.Cells.Item(455,8).Value = 9
$ws1.Cells.Item(455,9).Value = 4

$ws1.Cells.Item(456,1).Value = "ABC"
$ws1.Cells.Item(456,2).Value = 1
$ws1.Cells.Item(456,3).Value = 1
$ws1.Cells.Item(456,4).Value = 1
$ws1.Cells.Item(456,5).Value = 1
$ws1.Cells.Item(456,6).Value = 1
$ws1.Cells.Item(456,7).Value = 1
$ws1.Cells.Item(456,8).Value = 4
$ws1.Cells.Item(456,9).Value = 5

$ws1.Cells.Item(457,1).Value = "Ituano"
$ws1.Cells.Item(457,2).Value = 0
$ws1.Cells.Item(457,3).Value = 1
$ws1.Cells.Item(457,4).Value = 1
$ws1.Cells.Item(457,5).Value = 1
$ws1.Cells.Item(457,6).Value = 1
$ws1.Cells.Item(457,7).Value = 1
$ws1.Cells.Item(457,8).Value = 5
$ws1.Cells.Item(457,9).Value = 4

$ws1.Cells.Item(458,1).Value = "Londrina"
$ws1.Cells.Item(458,2).Value = 1
$ws1.Cells.Item(458,3).Value = 0
$ws1.Cells.Item(458,4).Value = 1
$ws1.Cells.Item(458,5).Value = 0
$ws1.Cells.Item(458,6).Value = 0
$ws1.Cells.Item(458,7).Value = 2
$ws1.Cells.Item(458,8).Value = 5
$ws1.Cells.Item(458,9).Value = 4

$ws1.Cells.Item(459,1).Value = "Atletico-GO"
$ws1.Cells.Item(459,2).Value = 0
$ws1.Cells.Item(459,3).Value = 1
$ws1.Cells.Item(459,4).Value = 0
$ws1.Cells.Item(459,5).Value = 0
$ws1.Cells.Item(459,6).Value = 2
$ws1.Cells.Item(459,7).Value = 0
$ws1.Cells.Item(459,8).Value = 4
$ws1.Cells.Item(459,9).Value = 5

$ws1.Cells.Item(460,1).Value = "Vila Nova"
$ws1.Cells.Item(460,2).Value = 1
$ws1.Cells.Item(460,3).Value = 1
$ws1.Cells.Item(460,4).Value = 0
$ws1.Cells.Item(460,5).Value = 0
$ws1.Cells.Item(460,6).Value = 2
$ws1.Cells.Item(460,7).Value = 0
$ws1.Cells.Item(460,8).Value = 3
$ws1.Cells.Item(460,9).Value = 1

$ws1.Cells.Item(461,1).Value = "Mirassol"
$ws1.Cells.Item(461,2).Value = 0
$ws1.Cells.Item(461,3).Value = 0
$ws1.Cells.Item(461,4).Value = 1
$ws1.Cells.Item(461,5).Value = 0
$ws1.Cells.Item(461,6).Value = 0
$ws1.Cells.Item(461,7).Value = 2
$ws1.Cells.Item(461,8).Value = 1
$ws1.Cells.Item(461,9).Value = 3

$ws1.Cells.Item(462,1).Value = "Vitoria"
$ws1.Cells.Item(462,2).Value = 1
$ws1.Cells.Item(462,3).Value = 1
$ws1.Cells.Item(462,4).Value = 0
$ws1.Cells.Item(462,5).Value = 0
$ws1.Cells.Item(462,6).Value = 2
$ws1.Cells.Item(462,7).Value = 0
$ws1.Cells.Item(462,8).Value = 11
$ws1.Cells.Item(462,9).Value = 1

$ws1.Cells.Item(463,1).Value = "Botafogo"
$ws1.Cells.Item(463,2).Value = 0
$ws1.Cells.Item(463,3).Value = 0
$ws1.Cells.Item(463,4).Value = 1
$ws1.Cells.Item(463,5).Value = 0
$ws1.Cells.Item(463,6).Value = 0
$ws1.Cells.Item(463,7).Value = 2
$ws1.Cells.Item(463,8).Value = 1
$ws1.Cells.Item(463,9).Value = 11

$ws1.Cells.Item(464,1).Value = "Guarani"
$ws1.Cells.Item(464,2).Value = 1
$ws1.Cells.Item(464,3).Value = 1
$ws1.Cells.Item(464,4).Value = 1
$ws1.Cells.Item(464,5).Value = 1
$ws1.Cells.Item(464,6).Value = 3
$ws1.Cells.Item(464,7).Value = 1
$ws1.Cells.Item(464,8).Value = 3
$ws1.Cells.Item(464,9).Value = 8

$ws1.Cells.Item(465,1).Value = "Sport"
$ws1.Cells.Item(465,2).Value = 0
$ws1.Cells.Item(465,3).Value = 1
$ws1.Cells.Item(465,4).Value = 1
$ws1.Cells.Item(465,5).Value = 1
$ws1.Cells.Item(465,6).Value = 1
$ws1.Cells.Item(465,7).Value = 3
$ws1.Cells.Item(465,8).Value = 8
$ws1.Cells.Item(465,9).Value = 3

$ws1.Cells.Item(466,1).Value = "Avai"
$ws1.Cells.Item(466,2).Value = 1
$ws1.Cells.Item(466,3).Value = 1
$ws1.Cells.Item(466,4).Value = 1
$ws1.Cells.Item(466,5).Value = 1
$ws1.Cells.Item(466,6).Value = 4
$ws1.Cells.Item(466,7).Value = 2
$ws1.Cells.Item(466,8).Value = 4
$ws1.Cells.Item(466,9).Value = 11

$ws1.Cells.Item(467,1).Value = "Tombense"
$ws1.Cells.Item(467,2).Value = 0
$ws1.Cells.Item(467,3).Value = 1
$ws1.Cells.Item(467,4).Value = 1
$ws1.Cells.Item(467,5).Value = 1
$ws1.Cells.Item(467,6).Value = 2
$ws1.Cells.Item(467,7).Value = 4
$ws1.Cells.Item(467,8).Value = 11
$ws1.Cells.Item(467,9).Value = 4

$ws1.Cells.Item(468,1).Value = "ABC"
$ws1.Cells.Item(468,2).Value = 1
$ws1.Cells.Item(468,3).Value = 1
$ws1.Cells.Item(468,4).Value = 1
$ws1.Cells.Item(468,5).Value = 1
$ws1.Cells.Item(468,6).Value = 1
$ws1.Cells.Item(468,7).Value = 2
$ws1.Cells.Item(468,8).Value = 8
$ws1.Cells.Item(468,9).Value = 3

$ws1.Cells.Item(469,1).Value = "CRB"
$ws1.Cells.Item(469,2).Value = 0
$ws1.Cells.Item(469,3).Value = 1
$ws1.Cells.Item(469,4).Value = 1
$ws1.Cells.Item(469,5).Value = 1
$ws1.Cells.Item(469,6).Value = 2
$ws1.Cells.Item(469,7).Value = 1
$ws1.Cells.Item(469,8).Value = 3
$ws1.Cells.Item(469,9).Value = 8

$ws1.Cells.Item(470,1).Value = "Ituano"
$ws1.Cells.Item(470,2).Value = 1
$ws1.Cells.Item(470,3).Value = 1
$ws1.Cells.Item(470,4).Value = 0
$ws1.Cells.Item(470,5).Value = 0
$ws1.Cells.Item(470,6).Value = 3
$ws1.Cells.Item(470,7).Value = 0
$ws1.Cells.Item(470,8).Value = 10
$ws1.Cells.Item(470,9).Value = 1

$ws1.Cells.Item(471,1).Value = "Criciuma"
$ws1.Cells.Item(471,2).Value = 0
$ws1.Cells.Item(471,3).Value = 0
$ws1.Cells.Item(471,4).Value = 1
$ws1.Cells.Item(471,5).Value = 0
$ws1.Cells.Item(471,6).Value = 0
$ws1.Cells.Item(471,7).Value = 3
$ws1.Cells.Item(471,8).Value = 1
$ws1.Cells.Item(471,9).Value = 10

$ws1.Cells.Item(472,1).Value = "Juventude"
$ws1.Cells.Item(472,2).Value = 1
$ws1.Cells.Item(472,3).Value = 0
$ws1.Cells.Item(472,4).Value = 0
$ws1.Cells.Item(472,5).Value = 0
$ws1.Cells.Item(472,6).Value = 0
$ws1.Cells.Item(472,7).Value = 0
$ws1.Cells.Item(472,8).Value = 10
$ws1.Cells.Item(472,9).Value = 2

$ws1.Cells.Item(473,1).Value = "Sampaio Correia"
$ws1.Cells.Item(473,2).Value = 0
$ws1.Cells.Item(473,3).Value = 0
$ws1.Cells.Item(473,4).Value = 0
$ws1.Cells.Item(473,5).Value = 0
$ws1.Cells.Item(473,6).Value = 0
$ws1.Cells.Item(473,7).Value = 0
$ws1.Cells.Item(473,8).Value = 2
$ws1.Cells.Item(473,9).Value = 10

$ws1.Cells.Item(474,1).Value = "Novohorizontino"
$ws1.Cells.Item(474,2).Value = 1
$ws1.Cells.Item(474,3).Value = 1
$ws1.Cells.Item(474,4).Value = 1
$ws1.Cells.Item(474,5).Value = 1
$ws1.Cells.Item(474,6).Value = 1
$ws1.Cells.Item(474,7).Value = 2
$ws1.Cells.Item(474,8).Value = 11
$ws1.Cells.Item(474,9).Value = 3

$ws1.Cells.Item(475,1).Value = "Chapecoense"
$ws1.Cells.Item(475,2).Value = 0
$ws1.Cells.Item(475,3).Value = 1
$ws1.Cells.Item(475,4).Value = 1
$ws1.Cells.Item(475,5).Value = 1
$ws1.Cells.Item(475,6).Value = 2
$ws1.Cells.Item(475,7).Value = 1
$ws1.Cells.Item(475,8).Value = 3
$ws1.Cells.Item(475,9).Value = 11

$ws1.Cells.Item(476,1).Value = "Ceara"
$ws1.Cells.Item(476,2).Value = 1
$ws1.Cells.Item(476,3).Value = 1
$ws1.Cells.Item(476,4).Value = 1
$ws1.Cells.Item(476,5).Value = 1
$ws1.Cells.Item(476,6).Value = 1
$ws1.Cells.Item(476,7).Value = 1
$ws1.Cells.Item(476,8).Value = 4
$ws1.Cells.Item(476,9).Value = 5

$ws1.Cells.Item(477,1).Value = "Ponte Preta"
$ws1.Cells.Item(477,2).Value = 0
$ws1.Cells.Item(477,3).Value = 1
$ws1.Cells.Item(477,4).Value = 1
$ws1.Cells.Item(477,5).Value = 1
$ws1.Cells.Item(477,6).Value = 1
$ws1.Cells.Item(477,7).Value = 1
$ws1.Cells.Item(477,8).Value = 5
$ws1.Cells.Item(477,9).Value = 4

$ws1.Cells.Item(478,1).Value = "Mirassol"
$ws1.Cells.Item(478,2).Value = 1
$ws1.Cells.Item(478,3).Value = 0
$ws1.Cells.Item(478,4).Value = 1
$ws1.Cells.Item(478,5).Value = 0
$ws1.Cells.Item(478,6).Value = 0
$ws1.Cells.Item(478,7).Value = 1
$ws1.Cells.Item(478,8).Value = 5
$ws1.Cells.Item(478,9).Value = 5

$ws1.Cells.Item(479,1).Value = "Juventude"
$ws1.Cells.Item(479,2).Value = 0
$ws1.Cells.Item(479,3).Value = 1
$ws1.Cells.Item(479,4).Value = 0
$ws1.Cells.Item(479,5).Value = 0
$ws1.Cells.Item(479,6).Value = 1
$ws1.Cells.Item(479,7).Value = 0
$ws1.Cells.Item(479,8).Value = 5
$ws1.Cells.Item(479,9).Value = 5

$ws1.Cells.Item(480,1).Value = "Criciuma"
$ws1.Cells.Item(480,2).Value = 1
$ws1.Cells.Item(480,3).Value = 1
$ws1.Cells.Item(480,4).Value = 0
$ws1.Cells.Item(480,5).Value = 0
$ws1.Cells.Item(480,6).Value = 1
$ws1.Cells.Item(480,7).Value = 0
$ws1.Cells.Item(480,8).Value = 5
$ws1.Cells.Item(480,9).Value = 2

$ws1.Cells.Item(481,1).Value = "Vila Nova"
$ws1.Cells.Item(481,2).Value = 0
$ws1.Cells.Item(481,3).Value = 0
$ws1.Cells.Item(481,4).Value = 1
$ws1.Cells.Item(481,5).Value = 0
$ws1.Cells.Item(481,6).Value = 0
$ws1.Cells.Item(481,7).Value = 1
$ws1.Cells.Item(481,8).Value = 2
$ws1.Cells.Item(481,9).Value = 5

$ws1.Cells.Item(482,1).Value = "Ponte Preta"
$ws1.Cells.Item(482,2).Value = 1
$ws1.Cells.Item(482,3).Value = 1
$ws1.Cells.Item(482,4).Value = 0
$ws1.Cells.Item(482,5).Value = 0
$ws1.Cells.Item(482,6).Value = 1
$ws1.Cells.Item(482,7).Value = 0
$ws1.Cells.Item(482,8).Value = 2
$ws1.Cells.Item(482,9).Value = 13

$ws1.Cells.Item(483,1).Value = "Londrina"
$ws1.Cells.Item(483,2).Value = 0
$ws1.Cells.Item(483,3).Value = 0
$ws1.Cells.Item(483,4).Value = 1
$ws1.Cells.Item(483,5).Value = 0
$ws1.Cells.Item(483,6).Value = 0
$ws1.Cells.Item(483,7).Value = 1
$ws1.Cells.Item(483,8).Value = 13
$ws1.Cells.Item(483,9).Value = 2

$ws1.Cells.Item(484,1).Value = "Sport"
$ws1.Cells.Item(484,2).Value = 1
$ws1.Cells.Item(484,3).Value = 1
$ws1.Cells.Item(484,4).Value = 1
$ws1.Cells.Item(484,5).Value = 1
$ws1.Cells.Item(484,6).Value = 1
$ws1.Cells.Item(484,7).Value = 2
$ws1.Cells.Item(484,8).Value = 11
$ws1.Cells.Item(484,9).Value = 3

$ws1.Cells.Item(485,1).Value = "Ituano"
$ws1.Cells.Item(485,2).Value = 0
$ws1.Cells.Item(485,3).Value = 1
$ws1.Cells.Item(485,4).Value = 1
$ws1.Cells.Item(485,5).Value = 1
$ws1.Cells.Item(485,6).Value = 2
$ws1.Cells.Item(485,7).Value = 1
$ws1.Cells.Item(485,8).Value = 3
$ws1.Cells.Item(485,9).Value = 11

$ws1.Cells.Item(486,1).Value = "Sampaio Correia"
$ws1.Cells.Item(486,2).Value = 1
$ws1.Cells.Item(486,3).Value = 1
$ws1.Cells.Item(486,4).Value = 1
$ws1.Cells.Item(486,5).Value = 1
$ws1.Cells.Item(486,6).Value = 1
$ws1.Cells.Item(486,7).Value = 1
$ws1.Cells.Item(486,8).Value = 9
$ws1.Cells.Item(486,9).Value = 4

$ws1.Cells.Item(487,1).Value = "Guarani"
$ws1.Cells.Item(487,2).Value = 0
$ws1.Cells.Item(487,3).Value = 1
$ws1.Cells.Item(487,4).Value = 1
$ws1.Cells.Item(487,5).Value = 1
$ws1.Cells.Item(487,6).Value = 1
$ws1.Cells.Item(487,7).Value = 1
$ws1.Cells.Item(487,8).Value = 4
$ws1.Cells.Item(487,9).Value = 9

$ws1.Cells.Item(488,1).Value = "Tombense"
$ws1.Cells.Item(488,2).Value = 1
$ws1.Cells.Item(488,3).Value = 1
$ws1.Cells.Item(488,4).Value = 1
$ws1.Cells.Item(488,5).Value = 1
$ws1.Cells.Item(488,6).Value = 2
$ws1.Cells.Item(488,7).Value = 2
$ws1.Cells.Item(488,8).Value = 10
$ws1.Cells.Item(488,9).Value = 4

$ws1.Cells.Item(489,1).Value = "Ceara"
$ws1.Cells.Item(489,2).Value = 0
$ws1.Cells.Item(489,3).Value = 1
$ws1.Cells.Item(489,4).Value = 1
$ws1.Cells.Item(489,5).Value = 1
$ws1.Cells.Item(489,6).Value = 2
$ws1.Cells.Item(489,7).Value = 2
$ws1.Cells.Item(489,8).Value = 4
$ws1.Cells.Item(489,9).Value = 10

$ws1.Cells.Item(490,1).Value = "CRB"
$ws1.Cells.Item(490,2).Value = 1
$ws1.Cells.Item(490,3).Value = 1
$ws1.Cells.Item(490,4).Value = 0
$ws1.Cells.Item(490,5).Value = 0
$ws1.Cells.Item(490,6).Value = 1
$ws1.Cells.Item(490,7).Value = 0
$ws1.Cells.Item(490,8).Value = 5
$ws1.Cells.Item(490,9).Value = 3

$ws1.Cells.Item(491,1).Value = "Novohorizontino"
$ws1.Cells.Item(491,2).Value = 0
$ws1.Cells.Item(491,3).Value = 0
$ws1.Cells.Item(491,4).Value = 1
$ws1.Cells.Item(491,5).Value = 0
$ws1.Cells.Item(491,6).Value = 0
$ws1.Cells.Item(491,7).Value = 1
$ws1.Cells.Item(491,8).Value = 3
$ws1.Cells.Item(491,9).Value = 5
$ws2 = $wb.Worksheets.Item("BD_Jogo")

$ws2.Cells.Item(221,1).Value = 0
$ws2.Cells.Item(221,2).Value = 0
$ws2.Cells.Item(221,3).Value = 7
$ws2.Cells.Item(221,4).Value = "Tombense"
$ws2.Cells.Item(221,5).Value = "Sport"

$ws2.Cells.Item(222,1).Value = 0
$ws2.Cells.Item(222,2).Value = 2
$ws2.Cells.Item(222,3).Value = 10
$ws2.Cells.Item(222,4).Value = "Criciuma"
$ws2.Cells.Item(222,5).Value = "Londrina"

$ws2.Cells.Item(223,1).Value = 0
$ws2.Cells.Item(223,2).Value = 0
$ws2.Cells.Item(223,3).Value = 6
$ws2.Cells.Item(223,4).Value = "Botafogo"
$ws2.Cells.Item(223,5).Value = "Ponte Preta"

$ws2.Cells.Item(224,1).Value = 1
$ws2.Cells.Item(224,2).Value = 2
$ws2.Cells.Item(224,3).Value = 14
$ws2.Cells.Item(224,4).Value = "Sampaio Correia"
$ws2.Cells.Item(224,5).Value = "CRB"

$ws2.Cells.Item(225,1).Value = 1
$ws2.Cells.Item(225,2).Value = 2
$ws2.Cells.Item(225,3).Value = 8
$ws2.Cells.Item(225,4).Value = "Vila Nova"
$ws2.Cells.Item(225,5).Value = "Avai"

$ws2.Cells.Item(226,1).Value = 0
$ws2.Cells.Item(226,2).Value = 1
$ws2.Cells.Item(226,3).Value = 8
$ws2.Cells.Item(226,4).Value = "Chapecoense"
$ws2.Cells.Item(226,5).Value = "Atletico-GO"

$ws2.Cells.Item(227,1).Value = 0
$ws2.Cells.Item(227,2).Value = 1
$ws2.Cells.Item(227,3).Value = 2
$ws2.Cells.Item(227,4).Value = "Guarani"
$ws2.Cells.Item(227,5).Value = "Juventude"

$ws2.Cells.Item(228,1).Value = 0
$ws2.Cells.Item(228,2).Value = 1
$ws2.Cells.Item(228,3).Value = 13
$ws2.Cells.Item(228,4).Value = "Vitoria"
$ws2.Cells.Item(228,5).Value = "Ceara"

$ws2.Cells.Item(229,1).Value = 1
$ws2.Cells.Item(229,2).Value = 2
$ws2.Cells.Item(229,3).Value = 9
$ws2.Cells.Item(229,4).Value = "ABC"
$ws2.Cells.Item(229,5).Value = "Ituano"

$ws2.Cells.Item(230,1).Value = 0
$ws2.Cells.Item(230,2).Value = 2
$ws2.Cells.Item(230,3).Value = 9
$ws2.Cells.Item(230,4).Value = "Londrina"
$ws2.Cells.Item(230,5).Value = "Atletico-GO"

$ws2.Cells.Item(231,1).Value = 0
$ws2.Cells.Item(231,2).Value = 2
$ws2.Cells.Item(231,3).Value = 4
$ws2.Cells.Item(231,4).Value = "Vila Nova"
$ws2.Cells.Item(231,5).Value = "Mirassol"

$ws2.Cells.Item(232,1).Value = 0
$ws2.Cells.Item(232,2).Value = 2
$ws2.Cells.Item(232,3).Value = 12
$ws2.Cells.Item(232,4).Value = "Vitoria"
$ws2.Cells.Item(232,5).Value = "Botafogo"

$ws2.Cells.Item(233,1).Value = 1
$ws2.Cells.Item(233,2).Value = 4
$ws2.Cells.Item(233,3).Value = 11
$ws2.Cells.Item(233,4).Value = "Guarani"
$ws2.Cells.Item(233,5).Value = "Sport"

$ws2.Cells.Item(234,1).Value = 1
$ws2.Cells.Item(234,2).Value = 6
$ws2.Cells.Item(234,3).Value = 15
$ws2.Cells.Item(234,4).Value = "Avai"
$ws2.Cells.Item(234,5).Value = "Tombense"

$ws2.Cells.Item(235,1).Value = 1
$ws2.Cells.Item(235,2).Value = 3
$ws2.Cells.Item(235,3).Value = 11
$ws2.Cells.Item(235,4).Value = "ABC"
$ws2.Cells.Item(235,5).Value = "CRB"

$ws2.Cells.Item(236,1).Value = 0
$ws2.Cells.Item(236,2).Value = 3
$ws2.Cells.Item(236,3).Value = 11
$ws2.Cells.Item(236,4).Value = "Ituano"
$ws2.Cells.Item(236,5).Value = "Criciuma"

$ws2.Cells.Item(237,1).Value = 0
$ws2.Cells.Item(237,2).Value = 0
$ws2.Cells.Item(237,3).Value = 12
$ws2.Cells.Item(237,4).Value = "Juventude"
$ws2.Cells.Item(237,5).Value = "Sampaio Correia"

$ws2.Cells.Item(238,1).Value = 1
$ws2.Cells.Item(238,2).Value = 3
$ws2.Cells.Item(238,3).Value = 14
$ws2.Cells.Item(238,4).Value = "Novohorizontino"
$ws2.Cells.Item(238,5).Value = "Chapecoense"

$ws2.Cells.Item(239,1).Value = 1
$ws2.Cells.Item(239,2).Value = 2
$ws2.Cells.Item(239,3).Value = 9
$ws2.Cells.Item(239,4).Value = "Ceara"
$ws2.Cells.Item(239,5).Value = "Ponte Preta"

$ws2.Cells.Item(240,1).Value = 0
$ws2.Cells.Item(240,2).Value = 1
$ws2.Cells.Item(240,3).Value = 10
$ws2.Cells.Item(240,4).Value = "Mirassol"
$ws2.Cells.Item(240,5).Value = "Juventude"

$ws2.Cells.Item(241,1).Value = 0
$ws2.Cells.Item(241,2).Value = 1
$ws2.Cells.Item(241,3).Value = 7
$ws2.Cells.Item(241,4).Value = "Criciuma"
$ws2.Cells.Item(241,5).Value = "Vila Nova"

$ws2.Cells.Item(242,1).Value = 0
$ws2.Cells.Item(242,2).Value = 1
$ws2.Cells.Item(242,3).Value = 15
$ws2.Cells.Item(242,4).Value = "Ponte Preta"
$ws2.Cells.Item(242,5).Value = "Londrina"

$ws2.Cells.Item(243,1).Value = 1
$ws2.Cells.Item(243,2).Value = 3
$ws2.Cells.Item(243,3).Value = 14
$ws2.Cells.Item(243,4).Value = "Sport"
$ws2.Cells.Item(243,5).Value = "Ituano"

$ws2.Cells.Item(244,1).Value = 1
$ws2.Cells.Item(244,2).Value = 2
$ws2.Cells.Item(244,3).Value = 13
$ws2.Cells.Item(244,4).Value = "Sampaio Correia"
$ws2.Cells.Item(244,5).Value = "Guarani"

$ws2.Cells.Item(245,1).Value = 1
$ws2.Cells.Item(245,2).Value = 4
$ws2.Cells.Item(245,3).Value = 14
$ws2.Cells.Item(245,4).Value = "Tombense"
$ws2.Cells.Item(245,5).Value = "Ceara"

$ws2.Cells.Item(246,1).Value = 0
$ws2.Cells.Item(246,2).Value = 1
$ws2.Cells.Item(246,3).Value = 8
$ws2.Cells.Item(246,4).Value = "CRB"
$ws2.Cells.Item(246,5).Value = "Novohorizontino"
